{"js": "/* global Word */\n// Replaces the paragraph \"Hola que tal \" (a placeholder list item) with five\n// fully-authored list paragraphs describing the templates/CRUD workflow,\n// matching the target revision of Documentacion.docx. The replacement is\n// performed by injecting a literal OOXML fragment so the exact run/proofErr/\n// bookmark structure of the target document is reproduced faithfully.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that still contains the old placeholder text.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Hola que tal\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find the paragraph containing \"Hola que tal\".');\n}\n\nconst FRAGMENT = `<w:p><w:pPr><w:pStyle w:val=\"Default\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"8\"/></w:numPr><w:spacing w:after=\"160\"/><w:jc w:val=\"both\"/><w:rPr><w:color w:val=\"auto\"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t xml:space=\"preserve\">Tras crear esa primera aproximaci\u00f3n a la </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t>pagina</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t xml:space=\"preserve\"> de </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t>Login</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t xml:space=\"preserve\"> crearemos dentro de nuestros \u2018</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t>templates</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t xml:space=\"preserve\">\u2019, las plantillas correspondientes para el CRUD de los diferentes tipos de modelos en cada app. Para ello, primero vamos a hacer uso de un \u2018Factory </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t>faker</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t>\u2019 que, tras instalarlo, nos va a poblar nuestros modelos para empezar a hacer pruebas.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Default\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"8\"/></w:numPr><w:spacing w:after=\"160\"/><w:jc w:val=\"both\"/><w:rPr><w:color w:val=\"auto\"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t>Cabe mencionar que la inserci\u00f3n de los formularios lo hacemos gracias a la herramienta de Django que nos permite simplemente insertar el formulario en la plantilla HTML con la instrucci\u00f3n \u2018{</w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t xml:space=\"preserve\">{ </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t>form</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:proofErr w:type=\"gramEnd\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t xml:space=\"preserve\"> }}\u2019.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Default\"/><w:spacing w:after=\"160\"/><w:jc w:val=\"both\"/><w:rPr><w:color w:val=\"auto\"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val=\"Default\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"8\"/></w:numPr><w:spacing w:after=\"160\"/><w:jc w:val=\"both\"/><w:rPr><w:color w:val=\"auto\"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space=\"preserve\">Despu\u00e9s de crear todos los </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t>templates</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t xml:space=\"preserve\"> de los correspondientes </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t>CURDs</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t xml:space=\"preserve\"> creamos un simple primer men\u00fa para cada uno de ellos.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Default\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"8\"/></w:numPr><w:spacing w:after=\"160\"/><w:jc w:val=\"both\"/><w:rPr><w:color w:val=\"auto\"/></w:rPr></w:pPr><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>`;\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>${FRAGMENT}<w:sectPr/></w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\ntarget.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Replaces the paragraph \"Hola que tal \" (a placeholder list item) with five\n# fully-authored list paragraphs describing the templates/CRUD workflow,\n# matching the target revision of Documentacion.docx. The replacement is\n# performed by locating the paragraph via Find, expanding to the full\n# paragraph range, and injecting a literal OOXML fragment via InsertXML so\n# the exact run/proofErr/bookmark structure of the target document is\n# reproduced faithfully.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Hola que tal\"\n$find.Forward = $true\n$find.Wrap = 1\n$found = $find.Execute()\nif (-not $found) {\n    throw 'Could not find the paragraph containing \"Hola que tal\".'\n}\n\n$targetRange = $find.Parent.Duplicate\n[void]$targetRange.Expand(4)\n\n$fragment = @\"\n<w:p><w:pPr><w:pStyle w:val=\"Default\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"8\"/></w:numPr><w:spacing w:after=\"160\"/><w:jc w:val=\"both\"/><w:rPr><w:color w:val=\"auto\"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t xml:space=\"preserve\">Tras crear esa primera aproximaci\u00f3n a la </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t>pagina</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t xml:space=\"preserve\"> de </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t>Login</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t xml:space=\"preserve\"> crearemos dentro de nuestros \u2018</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t>templates</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t xml:space=\"preserve\">\u2019, las plantillas correspondientes para el CRUD de los diferentes tipos de modelos en cada app. Para ello, primero vamos a hacer uso de un \u2018Factory </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t>faker</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t>\u2019 que, tras instalarlo, nos va a poblar nuestros modelos para empezar a hacer pruebas.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Default\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"8\"/></w:numPr><w:spacing w:after=\"160\"/><w:jc w:val=\"both\"/><w:rPr><w:color w:val=\"auto\"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t>Cabe mencionar que la inserci\u00f3n de los formularios lo hacemos gracias a la herramienta de Django que nos permite simplemente insertar el formulario en la plantilla HTML con la instrucci\u00f3n \u2018{</w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t xml:space=\"preserve\">{ </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t>form</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:proofErr w:type=\"gramEnd\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t xml:space=\"preserve\"> }}\u2019.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Default\"/><w:spacing w:after=\"160\"/><w:jc w:val=\"both\"/><w:rPr><w:color w:val=\"auto\"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val=\"Default\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"8\"/></w:numPr><w:spacing w:after=\"160\"/><w:jc w:val=\"both\"/><w:rPr><w:color w:val=\"auto\"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space=\"preserve\">Despu\u00e9s de crear todos los </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t>templates</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t xml:space=\"preserve\"> de los correspondientes </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t>CURDs</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t xml:space=\"preserve\"> creamos un simple primer men\u00fa para cada uno de ellos.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Default\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"8\"/></w:numPr><w:spacing w:after=\"160\"/><w:jc w:val=\"both\"/><w:rPr><w:color w:val=\"auto\"/></w:rPr></w:pPr><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>\n\"@\n\n$ooxml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>$fragment<w:sectPr/></w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n\"@\n\n$targetRange.InsertXML($ooxml)\n"}
